# Reaction-sensitivity results were regenerated with a fixed workflow.
# For both sheets ("NBR" and "BAR") the first 4 data rows (old cutoffs 1-4)
# are dropped, the remaining rows shift up, and the "Cutoff" index column
# (A) is renumbered starting at 0 while "Reaction_number" (C) gets the
# freshly computed values from the corrected pipeline.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NBR")
$ws.Rows("2:5").Delete()
$data = @(
    @(2, 0, 5, 851),
    @(3, 1, 6, 850),
    @(4, 2, 7, 850),
    @(5, 3, 8, 0),
    @(6, 4, 9, 845),
    @(7, 5, 10, 845),
    @(8, 6, 11, 843),
    @(9, 7, 12, 834),
    @(10, 8, 13, 832),
    @(11, 9, 14, 832),
    @(12, 10, 15, 833),
    @(13, 11, 16, 831),
    @(14, 12, 17, 831),
    @(15, 13, 18, 827),
    @(16, 14, 19, 826)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws = $wb.Worksheets.Item("BAR")
$ws.Rows("2:5").Delete()
$data = @(
    @(2, 0, 5, 689),
    @(3, 1, 6, 685),
    @(4, 2, 7, 683),
    @(5, 3, 8, 0),
    @(6, 4, 9, 681),
    @(7, 5, 10, 678),
    @(8, 6, 11, 679),
    @(9, 7, 12, 691),
    @(10, 8, 13, 693),
    @(11, 9, 14, 688),
    @(12, 10, 15, 685),
    @(13, 11, 16, 685),
    @(14, 12, 17, 687),
    @(15, 13, 18, 687),
    @(16, 14, 19, 684)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

